# Auto-generated: apply 2023-12-19 daily crime data update to violent-crime-full-year.xlsx
# Each sheet maps cell references (column J = year 2023 cumulative totals) to their new values.
$wb = $excel.ActiveWorkbook

$changes = @{
    'Citywide Totals' = @{ "J2" = 7413; "J3" = 7807; "J4" = 1697; "J5" = 612; "J6" = 10651; "J7" = 28180 }
    'By Neighborhood' = @{ "J2" = 222; "J8" = 1773; "J11" = 504; "J14" = 149; "J15" = 348; "J18" = 228; "J19" = 815; "J20" = 609; "J23" = 258; "J24" = 96; "J25" = 144; "J29" = 1502; "J30" = 98; "J31" = 298; "J33" = 1278; "J34" = 131; "J36" = 383; "J37" = 863; "J42" = 1198; "J47" = 205; "J48" = 315; "J50" = 170; "J51" = 355; "J54" = 558; "J55" = 443; "J57" = 134; "J60" = 166; "J63" = 82; "J64" = 189; "J65" = 708; "J68" = 60; "J73" = 275; "J76" = 400; "J77" = 198; "J79" = 770; "J83" = 566; "J84" = 234; "J85" = 1160; "J89" = 353; "J91" = 324; "J92" = 92; "J94" = 314; "J96" = 318; "J100" = 49; "J101" = 28180 }
    'Bridgeport' = @{ "J4" = 9; "J7" = 149 }
    'West Ridge' = @{ "J6" = 123; "J7" = 318 }
    'Belmont Cragin' = @{ "J3" = 87; "J6" = 240; "J7" = 504 }
    'Uptown' = @{ "J3" = 100; "J6" = 109; "J7" = 353 }
    'South Shore' = @{ "J2" = 308; "J3" = 419; "J6" = 332; "J7" = 1160 }
    'Little Village' = @{ "J3" = 201; "J6" = 307 }
    'Austin' = @{ "J2" = 466; "J5" = 45; "J6" = 657; "J7" = 1773 }
    'South Chicago' = @{ "J2" = 167; "J7" = 566 }
    'Garfield Park' = @{ "J3" = 426; "J6" = 454; "J7" = 1278 }
    'Grand Crossing' = @{ "J3" = 288; "J6" = 252; "J7" = 863 }
    'New City' = @{ "J3" = 189; "J7" = 708 }
    'Fuller Park' = @{ "J2" = 35; "J7" = 98 }
    'Gage Park' = @{ "J2" = 97; "J3" = 73; "J6" = 107; "J7" = 298 }
    'South Deering' = @{ "J2" = 71; "J7" = 234 }
    'Loop' = @{ "J6" = 258; "J7" = 558 }
    'Englewood' = @{ "J2" = 458; "J3" = 528; "J4" = 81; "J5" = 56; "J7" = 1502 }
    'Lake View' = @{ "J4" = 49; "J7" = 315 }
    'Chatham' = @{ "J2" = 200; "J3" = 231; "J6" = 316; "J7" = 815 }
    'River North' = @{ "J3" = 89; "J7" = 400 }
    'Humboldt Park' = @{ "J2" = 250; "J3" = 241; "J6" = 634; "J7" = 1198 }
    'Rogers Park' = @{ "J2" = 86 }
    'Lower West Side' = @{ "J7" = 443 }
    'Dunning' = @{ "J3" = 25; "J7" = 96 }
    'Douglas' = @{ "J3" = 85; "J7" = 258 }
    'Washington Park' = @{ "J2" = 85; "J3" = 132; "J7" = 324 }
    'Roseland' = @{ "J2" = 219; "J6" = 231; "J7" = 770 }
    'Near South Side' = @{ "J2" = 51; "J7" = 189 }
    'Chicago Lawn' = @{ "J2" = 169; "J3" = 199; "J5" = 18; "J7" = 609 }
    'Calumet Heights' = @{ "J4" = 12; "J7" = 228 }
    'Grand Boulevard' = @{ "J3" = 124; "J7" = 383 }
    'Wrigleyville' = @{ "J4" = 3; "J7" = 49 }
    'Garfield Ridge' = @{ "J3" = 35; "J7" = 131 }
    'West Loop' = @{ "J6" = 166; "J7" = 314 }
    'East Side' = @{ "J2" = 58; "J7" = 144 }
    'Kenwood' = @{ "J2" = 48; "J7" = 205 }
    'Brighton Park' = @{ "J6" = 161; "J7" = 348 }
    'Lincoln Square' = @{ "J2" = 43; "J4" = 26; "J7" = 170 }
    'Portage Park' = @{ "J3" = 71; "J7" = 275 }
    'Albany Park' = @{ "J3" = 53; "J7" = 222 }
    'West Elsdon' = @{ "J6" = 34; "J7" = 92 }
    'Little Italy, UIC' = @{ "J4" = 33; "J5" = 9; "J6" = 147; "J7" = 355 }
    'North Park' = @{ "J6" = 15; "J7" = 60 }
    'Mckinley Park' = @{ "J6" = 61; "J7" = 134 }
    'Morgan Park' = @{ "J2" = 57; "J7" = 166 }
    'Riverdale' = @{ "J2" = 74; "J7" = 198 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $changes[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}
